# Updates cryptos list prices / 1h-volume deltas (and a couple of swapped
# coin rows) per the scraped data refresh.
#
# Several Price cells hold short decimal-looking text (e.g. "0.580",
# "14.20") that must stay plain text, exactly as authored upstream
# (t="inlineStr"/shared-string, not numeric). Assigning such a value
# straight to .Value lets Excel's autodetect coerce it to a number and
# drop the trailing/insignificant digits (e.g. "0.580" -> 0.58). To keep
# it textual we assign with a leading apostrophe (Excel's literal-text
# escape) and then reset .Style to "Normal" so no stray NumberFormat /
# quote-prefix styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.648.43'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '2.243.38'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'321.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.47%  '
$ws.Range("D6").Value = "'101.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("D7").Value = "'0.580"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("D10").Value = "'37.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").Value = "'7.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("E13").Value = '  -2.16%  '
$ws.Range("D14").Value = '2.585.74'
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").Value = "'0.856"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.278.96'
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = "'14.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").Value = '43.535.43'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("D19").Value = "'13.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.38%  '
$ws.Range("E20").Value = '  +2.67%  '
$ws.Range("D21").Value = "'6.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").Value = "'65.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").Value = "'3.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").Value = "'236.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").Value = "'2.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = "'10.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("D28").Value = "'2.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.28%  '
$ws.Range("D29").Value = "'36.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.33%  '
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").Value = "'159.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.50%  '
$ws.Range("D32").Value = "'20.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").Value = "'0.0852"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.45%  '
$ws.Range("D34").Value = "'2.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.79%  '
$ws.Range("D35").Value = "'3.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.97%  '
$ws.Range("E36").Value = '  +8.79%  '
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("E39").Value = '  +1.92%  '
$ws.Range("D40").Value = "'4.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.64%  '
$ws.Range("D41").Value = "'15.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +20.34%  '
$ws.Range("D42").Value = "'0.0318"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.97%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").Value = '1.795.86'
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("E45").Value = '  -3.35%  '
$ws.Range("D46").Value = "'82.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.36%  '
$ws.Range("B47").Value = 'ordi'
$ws.Range("C47").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D47").Value = "'74.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.27%  '
$ws.Range("D48").Value = "'5.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.72%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = "'1.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.94%  '
$ws.Range("D50").Value = "'58.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.31%  '
$ws.Range("D51").Value = "'103.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.57%  '
